# Few changes on State and Visitor patterns
$p = $ppt.ActivePresentation

$guillemetOpen  = [char]0x2039
$guillemetClose = [char]0x203A

# --- Slide master 1 (Design 1): date placeholder + slide-number placeholder ---
$d1 = $p.Designs.Item(1).SlideMaster
$d1.Shapes.Item(2).TextFrame.TextRange.Text = "3/23/2022"
$d1.Shapes.Item(4).TextFrame.TextRange.Text = $guillemetOpen + "#" + $guillemetClose

# --- Slide master 2 (Design 2): date placeholder + slide-number placeholder ---
$d2 = $p.Designs.Item(2).SlideMaster
$d2.Shapes.Item(3).TextFrame.TextRange.Text = "3/23/2022"
$d2.Shapes.Item(4).TextFrame.TextRange.Text = $guillemetOpen + "#" + $guillemetClose

# --- Notes master: slide-number placeholder ---
$nm = $p.NotesMaster
$nm.Shapes.Item(6).TextFrame.TextRange.Text = $guillemetOpen + "#" + $guillemetClose

# --- Slide 16 ("Behavioral Patterns"): update languages on existing bullets and add "Visitor" ---
$s16 = $p.Slides.Item(16)
$bullets = $s16.Shapes.Item(2)
$tr = $bullets.TextFrame.TextRange

for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    $para.Font.LanguageID = 1033
}

$lastPara = $tr.Paragraphs($tr.Paragraphs().Count)
$lastPara.InsertAfter("`rVisitor")

Write-Host "edit complete"
